# Apply corrections to the "Before FS-DR" sheet per commit "code is update and fix"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Before FS-DR")

# --- Scalar value fixes (Data 1 / Data 7 blocks) ---
$ws.Range("B6").Value = 0.9789
$ws.Range("C6").Value = 0.9782999999999999

$ws.Range("B54").Value = 0.9781
$ws.Range("C54").Value = 0.9577

$ws.Range("B58").Value = 0.9772999999999999
$ws.Range("C58").Value = 0.9726

$ws.Range("B71").Value = 0.9804
$ws.Range("C71").Value = 0.9598

# --- SVM parameter (C=...) swaps ---
$ws.Range("L16").Value = "C=1"
$ws.Range("L21").Value = "C=10"
$ws.Range("L39").Value = "C=1"

# --- Remove stray duplicated "Data 8" block results (rows 75-84), keep fold labels in col A ---
$ws.Range("B75:P84").ClearContents()

# --- Remove stray duplicated "Data 9" block results (rows 87-96), keep fold labels in col A ---
$ws.Range("B87:P96").ClearContents()
